# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.030.30"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").Value = "1.834.30"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'243.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

# Row 6
$ws.Range("D6").Value = "'0.6292"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.16%  "

# Row 7
$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "'0.07614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.41%  "

# Row 9
$ws.Range("D9").Value = "'0.2935"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'22.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.24%  "

# Row 11
$ws.Range("D11").Value = "'0.07738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "

# Row 12
$ws.Range("D12").Value = "1.832.09"
$ws.Range("E12").Value = "  +0.23%  "

# Row 13
$ws.Range("D13").Value = "'4.963"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("D14").Value = "'0.6659"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001011"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.71%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'82.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17
$ws.Range("D17").Value = "'6.103"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "

# Row 18
$ws.Range("D18").Value = "29.026.87"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("D19").Value = "'226.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20
$ws.Range("D20").Value = "'12.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "

# Row 21
$ws.Range("D21").Value = "'0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").Value = "'7.196"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "

# Row 23
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "'159.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
$ws.Range("D25").Value = "'0.1406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.73%  "

# Row 26
$ws.Range("D26").Value = "'8.492"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "

# Row 27
$ws.Range("D27").Value = "'17.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "

# Row 28
$ws.Range("E28").Value = "  -0.22%  "

# Row 29
$ws.Range("D29").Value = "'4.101"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.66%  "

# Row 30
$ws.Range("D30").Value = "'4.017"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "

# Row 31
$ws.Range("D31").Value = "'1.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32
$ws.Range("D32").Value = "'0.05322"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$ws.Range("D33").Value = "'1.849"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("D34").Value = "'0.7379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("D35").Value = "'1.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "

# Row 36
$ws.Range("D36").Value = "'2.672"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "

# Row 37
$ws.Range("D37").Value = "1.238.64"
$ws.Range("E37").Value = "  -3.92%  "

# Row 38
$ws.Range("D38").Value = "'2.760"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.60%  "

# Row 39
$ws.Range("D39").Value = "'0.01784"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("D40").Value = "'6.351"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "

# Row 41
$ws.Range("D41").Value = "'0.8987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42
$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").Value = "'102.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

# Row 44
$ws.Range("D44").Value = "1.984.91"
$ws.Range("E44").Value = "  +0.42%  "

# Row 45
$ws.Range("E45").Value = "  +2.27%  "

# Row 46
$ws.Range("D46").Value = "'64.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("E47").Value = "  -0.15%  "

# Row 48
$ws.Range("D48").Value = "'0.4067"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.27%  "

# Row 49
$ws.Range("D49").Value = "'8.922"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.80%  "

# Row 50
$ws.Range("D50").Value = "'0.05784"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("D51").Value = "'6.727"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
